$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update odds in row 2 (Madura United vs Arema FC) ---
$ws.Cells.Item(2,7).Value = 2.4
$ws.Cells.Item(2,9).Value = 2.82
$ws.Cells.Item(2,11).Value = 2.02
$ws.Cells.Item(2,12).Value = 3.45
$ws.Cells.Item(2,15).Value = 1.27
$ws.Cells.Item(2,16).Value = 3.1
$ws.Cells.Item(2,18).Value = 1.85
$ws.Cells.Item(2,20).Value = 2.5
$ws.Cells.Item(2,24).Value = 12.5
$ws.Cells.Item(2,25).Value = 9
$ws.Cells.Item(2,26).Value = 26
$ws.Cells.Item(2,27).Value = 19
$ws.Cells.Item(2,29).Value = 9.75
$ws.Cells.Item(2,30).Value = 6.1
$ws.Cells.Item(2,34).Value = 9.25
$ws.Cells.Item(2,35).Value = 15
$ws.Cells.Item(2,37).Value = 35
$ws.Cells.Item(2,38).Value = 24
$ws.Cells.Item(2,40).Value = 4.3
$ws.Cells.Item(2,44).Value = 90
$ws.Cells.Item(2,46).Value = 2.47
$ws.Cells.Item(2,49).Value = 4.75
$ws.Cells.Item(2,50).Value = 15.5
$ws.Cells.Item(2,51).Value = 23
$ws.Cells.Item(2,52).Value = 75
$ws.Cells.Item(2,53).Value = 110

# --- Update odds in row 3 (Malut United vs Persis Solo) ---
$ws.Cells.Item(3,7).Value = 2.18
$ws.Cells.Item(3,8).Value = 2.92
$ws.Cells.Item(3,9).Value = 3.45
$ws.Cells.Item(3,10).Value = 2.67
$ws.Cells.Item(3,11).Value = 2.02
$ws.Cells.Item(3,12).Value = 3.95
$ws.Cells.Item(3,15).Value = 1.36
$ws.Cells.Item(3,16).Value = 2.67
$ws.Cells.Item(3,21).Value = 1.75
$ws.Cells.Item(3,22).Value = 1.85
$ws.Cells.Item(3,23).Value = 7.1
$ws.Cells.Item(3,24).Value = 10.5
$ws.Cells.Item(3,25).Value = 8.5
$ws.Cells.Item(3,26).Value = 22
$ws.Cells.Item(3,27).Value = 18
$ws.Cells.Item(3,28).Value = 29
$ws.Cells.Item(3,29).Value = 7.7
$ws.Cells.Item(3,30).Value = 5.7
$ws.Cells.Item(3,34).Value = 9
$ws.Cells.Item(3,35).Value = 18
$ws.Cells.Item(3,39).Value = 45
$ws.Cells.Item(3,40).Value = 4.05
$ws.Cells.Item(3,42).Value = 17.5
$ws.Cells.Item(3,43).Value = 40
$ws.Cells.Item(3,44).Value = 65
$ws.Cells.Item(3,45).Value = 200
$ws.Cells.Item(3,46).Value = 2.5
$ws.Cells.Item(3,47).Value = 6.5
$ws.Cells.Item(3,50).Value = 20
$ws.Cells.Item(3,51).Value = 25
$ws.Cells.Item(3,53).Value = 150

# --- Update odds in row 5 (Jeonnam vs Busan) ---
$ws.Cells.Item(5,17).Value = 1.85
$ws.Cells.Item(5,18).Value = 1.95

# --- Insert a new row 6 for the new match (Lunds vs Ostersund), shifting old rows 6-7 down to 7-8 ---
$ws.Rows.Item(6).Insert()

# --- Populate new row 6 ---
$ws.Cells.Item(6,1).Value = "0CEysv2s"
$ws.Cells.Item(6,2).Value = "21/11/2024"
$ws.Cells.Item(6,3).Value = "14:30"
$ws.Cells.Item(6,4).Value = "SWEDEN - SUPERETTAN"
$ws.Cells.Item(6,5).Value = "Lunds"
$ws.Cells.Item(6,6).Value = "Ostersund"
$ws.Cells.Item(6,7).Value = 3
$ws.Cells.Item(6,8).Value = 3.3
$ws.Cells.Item(6,9).Value = 2.2
$ws.Cells.Item(6,10).Value = 3.6
$ws.Cells.Item(6,11).Value = 2.2
$ws.Cells.Item(6,12).Value = 2.88
$ws.Cells.Item(6,13).Value = 1.05
$ws.Cells.Item(6,14).Value = 11
$ws.Cells.Item(6,15).Value = 1.29
$ws.Cells.Item(6,16).Value = 3.5
$ws.Cells.Item(6,17).Value = 1.93
$ws.Cells.Item(6,18).Value = 1.93
$ws.Cells.Item(6,19).Value = 1.36
$ws.Cells.Item(6,20).Value = 3
$ws.Cells.Item(6,21).Value = 1.67
$ws.Cells.Item(6,22).Value = 2.1
$ws.Cells.Item(6,23).Value = 10
$ws.Cells.Item(6,24).Value = 15
$ws.Cells.Item(6,25).Value = 11
$ws.Cells.Item(6,26).Value = 34
$ws.Cells.Item(6,27).Value = 23
$ws.Cells.Item(6,28).Value = 29
$ws.Cells.Item(6,29).Value = 11
$ws.Cells.Item(6,30).Value = 6.5
$ws.Cells.Item(6,31).Value = 13
$ws.Cells.Item(6,32).Value = 41
$ws.Cells.Item(6,33).Value = 151
$ws.Cells.Item(6,34).Value = 8.5
$ws.Cells.Item(6,35).Value = 12
$ws.Cells.Item(6,36).Value = 9.5
$ws.Cells.Item(6,37).Value = 21
$ws.Cells.Item(6,38).Value = 17
$ws.Cells.Item(6,39).Value = 26
$ws.Cells.Item(6,40).Value = 5
$ws.Cells.Item(6,41).Value = 17
$ws.Cells.Item(6,42).Value = 23
$ws.Cells.Item(6,43).Value = 51
$ws.Cells.Item(6,44).Value = 67
$ws.Cells.Item(6,45).Value = 151
$ws.Cells.Item(6,46).Value = 3
$ws.Cells.Item(6,47).Value = 7.5
$ws.Cells.Item(6,48).Value = 51
$ws.Cells.Item(6,49).Value = 4.33
$ws.Cells.Item(6,50).Value = 13
$ws.Cells.Item(6,51).Value = 21
$ws.Cells.Item(6,52).Value = 41
$ws.Cells.Item(6,53).Value = 51
$ws.Cells.Item(6,54).Value = 151
$ws.Cells.Item(6,55).Value = 126
$ws.Cells.Item(6,56).Value = 126
